# Make sure inventory UI matches with the real inventory
$wb = $excel.ActiveWorkbook

$wsNames = $wb.Worksheets.Item("NameTextEntities")
$wsDesc  = $wb.Worksheets.Item("DescriptionTextEntities")

# Fix the Japanese name for Puzzle Piece B (row 5, column D) - was incorrectly
# showing the same text as Puzzle Piece A.
$wsNames.Range("D5").Value = "パズルBなもの"

# Fix the Japanese descriptions for Puzzle Piece A and B (rows 4 and 5, column D)
# which previously shared the exact same (ambiguous) description text.
$wsDesc.Range("D4").Value = "Aパズルブロッカだ。パズル用みたいです。"
$wsDesc.Range("D5").Value = "Bパズルブロッカだ。パズル用みたいです。"

# Reflect the active-cell selections left by the editor in each sheet.
# Select the DescriptionTextEntities sheet's cell first, then the
# NameTextEntities sheet last so that NameTextEntities remains the active tab.
$wsDesc.Range("D4").Select() | Out-Null
$wsNames.Range("F7").Select() | Out-Null
